# Update cryptocurrency price and 1h volume change figures
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "48.304.08"
$ws.Range("E2").Value = "  +2.34%  "
$ws.Range("D3").Value = "2.515.83"
$ws.Range("E3").Value = "  +1.27%  "
$ws.Range("E4").Value = "  +0.03%  "
$origStyle = $ws.Range("D5").Style
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "321.79"
$ws.Range("D5").Style = $origStyle
$ws.Range("E5").Value = "  +0.32%  "
$origStyle = $ws.Range("D6").Style
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "108.89"
$ws.Range("D6").Style = $origStyle
$ws.Range("E6").Value = "  +0.77%  "
$origStyle = $ws.Range("D7").Style
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.529"
$ws.Range("D7").Style = $origStyle
$ws.Range("E8").Value = "  +0.05%  "
$ws.Range("E9").Value = "  +1.05%  "
$origStyle = $ws.Range("D10").Style
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "40.07"
$ws.Range("D10").Style = $origStyle
$ws.Range("E10").Value = "  +2.36%  "
$origStyle = $ws.Range("D11").Style
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.07"
$ws.Range("D11").Style = $origStyle
$ws.Range("E11").Value = "  +9.17%  "
$ws.Range("E12").Value = "  +1.30%  "
$ws.Range("E13").Value = "  +0.32%  "
$ws.Range("D15").Value = "2.909.04"
$ws.Range("E15").Value = "  +1.29%  "
$ws.Range("D16").Value = "2.517.81"
$ws.Range("E16").Value = "  +1.36%  "
$ws.Range("E17").Value = "  +0.49%  "
$ws.Range("D18").Value = "48.148.08"
$ws.Range("E18").Value = "  +2.22%  "
$origStyle = $ws.Range("D19").Style
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.15"
$ws.Range("D19").Style = $origStyle
$ws.Range("E19").Value = "  -1.94%  "
$origStyle = $ws.Range("D20").Style
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.79"
$ws.Range("D20").Style = $origStyle
$ws.Range("E20").Value = "  +2.75%  "
$ws.Range("D21").Value = "0.0₃0952"
$ws.Range("E21").Value = "  +0.98%  "
$origStyle = $ws.Range("D22").Style
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.75"
$ws.Range("D22").Style = $origStyle
$ws.Range("E22").Value = "  +0.28%  "
$origStyle = $ws.Range("D23").Style
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "72.31"
$ws.Range("D23").Style = $origStyle
$ws.Range("E23").Value = "  +2.56%  "
$origStyle = $ws.Range("D24").Style
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "277.97"
$ws.Range("D24").Style = $origStyle
$ws.Range("E24").Value = "  +13.14%  "
$origStyle = $ws.Range("D25").Style
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.57"
$ws.Range("D25").Style = $origStyle
$ws.Range("E25").Value = "  +0.88%  "
$origStyle = $ws.Range("D27").Style
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "25.93"
$ws.Range("D27").Style = $origStyle
$ws.Range("E27").Value = "  +0.97%  "
$ws.Range("E28").Value = "  +4.59%  "
$origStyle = $ws.Range("D29").Style
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.86"
$ws.Range("D29").Style = $origStyle
$ws.Range("E29").Value = "  -0.91%  "
$origStyle = $ws.Range("D30").Style
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "35.49"
$ws.Range("D30").Style = $origStyle
$ws.Range("E30").Value = "  +2.71%  "
$origStyle = $ws.Range("D31").Style
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.139"
$ws.Range("D31").Style = $origStyle
$ws.Range("E31").Value = "  -0.77%  "
$origStyle = $ws.Range("D32").Style
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "49.20"
$ws.Range("D32").Style = $origStyle
$ws.Range("E32").Value = "  -1.13%  "
$ws.Range("E33").Value = "  -3.69%  "
$origStyle = $ws.Range("D34").Style
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.38"
$ws.Range("D34").Style = $origStyle
$ws.Range("E34").Value = "  +0.95%  "
$ws.Range("E35").Value = "  +0.01%  "
$ws.Range("E36").Value = "  +0.63%  "
$ws.Range("E37").Value = "  +0.85%  "
$ws.Range("E38").Value = "  -2.14%  "
$ws.Range("E39").Value = "  +1.34%  "
$ws.Range("E40").Value = "  +0.21%  "
$origStyle = $ws.Range("D41").Style
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "122.27"
$ws.Range("D41").Style = $origStyle
$ws.Range("E41").Value = "  +3.59%  "
$origStyle = $ws.Range("D42").Style
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.21"
$ws.Range("D42").Style = $origStyle
$ws.Range("E42").Value = "  +0.36%  "
$origStyle = $ws.Range("D43").Style
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "21.44"
$ws.Range("D43").Style = $origStyle
$ws.Range("E43").Value = "  -6.96%  "
$ws.Range("E44").Value = "  +3.31%  "
$ws.Range("D45").Value = "2.002.74"
$ws.Range("E45").Value = "  +0.37%  "
$ws.Range("E46").Value = "  +4.85%  "
$ws.Range("E47").Value = "  +3.99%  "
$ws.Range("E48").Value = "  -1.19%  "
$ws.Range("E49").Value = "  -0.72%  "
$origStyle = $ws.Range("D50").Style
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "5.25"
$ws.Range("D50").Style = $origStyle
$ws.Range("E50").Value = "  +3.22%  "
$origStyle = $ws.Range("D51").Style
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "80.26"
$ws.Range("D51").Style = $origStyle
